$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data (columns A and B) between row 4 and row 8
$a4 = $ws.Range("A4").Value2
$b4 = $ws.Range("B4").Value2
$a8 = $ws.Range("A8").Value2
$b8 = $ws.Range("B8").Value2

$ws.Range("A4").Value = $a8
$ws.Range("B4").Value = $b8
$ws.Range("A8").Value = $a4
$ws.Range("B8").Value = $b4
